$d = $word.ActiveDocument

$replacements = @(
    @("89×12=", "21×26="),
    @("43×31=", "37×18="),
    @("13×85=", "45×15="),
    @("57×42=", "15×41="),
    @("73×37=", "16×34="),
    @("71×70=", "94×18="),
    @("52×49=", "78×74="),
    @("74×76=", "87×45="),
    @("76×87=", "63×20="),
    @("42×17=", "42×96="),
    @("13×34=", "98×67="),
    @("60×79=", "50×42="),
    @("29×39=", "80×19="),
    @("18×72=", "46×46="),
    @("46×12=", "11×47="),
    @("22×43=", "55×56="),
    @("61×98=", "26×23="),
    @("12×71=", "41×81="),
    @("29×57=", "67×23="),
    @("99×12=", "31×38="),
    @("35×25=", "85×98="),
    @("69×42=", "54×34="),
    @("82×98=", "26×23="),
    @("11×93=", "38×26="),
    @("48×44=", "39×29=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
